# Remove the "Implementierung - Datenfluss" slide from the deck.
$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "Implementierung - Datenfluss") {
                $targetIndex = $i
                break
            }
        }
    }
    if ($targetIndex -ne -1) {
        break
    }
}

if ($targetIndex -eq -1) {
    $targetIndex = 8
}

$p.Slides.Item($targetIndex).Delete()
